# Update data.xlsx from the QR tool output:
#  - overwrite row 2 with the latest scanned/demo record
#  - drop the old extra rows 3 and 4 (test rows) entirely
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing test rows first so row 2 becomes the last data row.
$ws.Rows("3:4").Delete()

# Refresh row 2 with the new record's values.
$ws.Range("A2").Value = "0kbwrl1cwnf8"
$ws.Range("B2").Value = "DEMOM9CC"
$ws.Range("C2").Value = "Hộ kinh doanh Trần Văn A"
$ws.Range("D2").Value = "02 Hòa Bình, Ninh Kiều, Cần Thơ"
$ws.Range("E2").Value = "https://www.google.com/maps/search/?api=1&query=02%20H%C3%B2a%20B%C3%ACnh%2C%20Ninh%20Ki%E1%BB%81u%2C%20C%E1%BA%A7n%20Th%C6%A1"
$ws.Range("F2").Value = "2025-08-15T01:18:36.878Z"
$ws.Range("G2").Value = "Dòng mẫu để thử"

# phone / cccd are numeric-looking strings that must stay text (leading
# zeros). A plain .Value assignment would auto-coerce them to numbers, so
# write them as a text formula in a scratch cell and paste-special the
# computed value back in - this keeps the literal text without touching
# the cell's number format/style.
$scratch = $ws.Range("ZZ1")

$scratch.Formula = '="0901234567"'
$scratch.Copy()
$ws.Range("H2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("I2").Value = "CN Cần Thơ II"

$scratch.Formula = '="012345678901"'
$scratch.Copy()
$ws.Range("J2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("K2").Value = "KH0001"
$ws.Range("L2").Value = "Demo User"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

$scratch.Clear()
$excel.CutCopyMode = $false
